$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 772.6
$ws.Range("I115").Value = 613.5
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 1840.5
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -273.5
$ws.Range("N115").Value = -12134

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3743.75
$ws.Range("I116").Value = 3043.4614
$ws.Range("J116").Value = 4571.364
$ws.Range("K116").Value = 3043.4614
$ws.Range("L116").Value = 4571.364
$ws.Range("M116").Value = 398.5385999999999
$ws.Range("N116").Value = -11455.364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1546.7333
$ws.Range("I45").Value = 1058.0834
$ws.Range("J45").Value = 3501.3333
$ws.Range("K45").Value = 1058.0834
$ws.Range("L45").Value = 3501.3333
$ws.Range("M45").Value = -681.0834
$ws.Range("N45").Value = -4255.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1174.1364
$ws.Range("I110").Value = 709.75
$ws.Range("J110").Value = 2412.5
$ws.Range("K110").Value = 709.75
$ws.Range("L110").Value = 2412.5
$ws.Range("M110").Value = 1335.25
$ws.Range("N110").Value = -6502.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 16951910
$ws.Range("I132").Value = 21279240
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 63837720
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -63835190
$ws.Range("N132").Value = -14660

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1687.6666
$ws.Range("I99").Value = 1060.8334
$ws.Range("M99").Value = 437.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1770.2307
$ws.Range("I107").Value = 540
$ws.Range("J107").Value = 2539.125
$ws.Range("K107").Value = 540
$ws.Range("L107").Value = 2539.125
$ws.Range("M107").Value = 1380
$ws.Range("N107").Value = -6379.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2870.0938
$ws.Range("I134").Value = 1793.8
$ws.Range("J134").Value = 6714
$ws.Range("K134").Value = 5381.4
$ws.Range("L134").Value = 20142
$ws.Range("M134").Value = -2846.4
$ws.Range("N134").Value = -25212

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 813400.6
$ws.Range("I6").Value = 1334667
$ws.Range("J6").Value = 31501
$ws.Range("K6").Value = 1334667
$ws.Range("L6").Value = 31501
$ws.Range("M6").Value = -1334554
$ws.Range("N6").Value = -31727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2428.4707
$ws.Range("I31").Value = 1496.5333
$ws.Range("J31").Value = 3759.8096
$ws.Range("K31").Value = 1496.5333
$ws.Range("L31").Value = 3759.8096
$ws.Range("M31").Value = -1201.5333
$ws.Range("N31").Value = -4349.809600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2428.4707
$ws.Range("I34").Value = 1496.5333
$ws.Range("J34").Value = 3759.8096
$ws.Range("K34").Value = 1496.5333
$ws.Range("L34").Value = 3759.8096
$ws.Range("M34").Value = -1294.5333
$ws.Range("N34").Value = -4163.809600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 17859512
$ws.Range("I58").Value = 1191.6666
$ws.Range("J58").Value = 31253252
$ws.Range("K58").Value = 1191.6666
$ws.Range("L58").Value = 31253252
$ws.Range("M58").Value = -988.6666
$ws.Range("N58").Value = -31253658

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2208.3845
$ws.Range("I134").Value = 1337.2273
$ws.Range("J134").Value = 6999.75
$ws.Range("K134").Value = 4011.6819
$ws.Range("L134").Value = 20999.25
$ws.Range("M134").Value = -1476.6819
$ws.Range("N134").Value = -26069.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 17859512
$ws.Range("I136").Value = 1191.6666
$ws.Range("J136").Value = 31253252
$ws.Range("K136").Value = 3574.9998
$ws.Range("L136").Value = 93759756
$ws.Range("M136").Value = -1024.9998
$ws.Range("N136").Value = -93764856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 9000
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = -6754
$ws.Range("N99").Value = -13492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1028.5745
$ws.Range("I107").Value = 638.26666
$ws.Range("J107").Value = 1211.5312
$ws.Range("K107").Value = 1914.79998
$ws.Range("L107").Value = 3634.5936
$ws.Range("M107").Value = 5.200019999999995
$ws.Range("N107").Value = -7474.5936

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 865.2174
$ws.Range("I114").Value = 597.875
$ws.Range("J114").Value = 1007.8
$ws.Range("K114").Value = 1793.625
$ws.Range("L114").Value = 3023.4
$ws.Range("M114").Value = 1460.375
$ws.Range("N114").Value = -9531.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 44945.582
$ws.Range("I129").Value = 5373.4
$ws.Range("J129").Value = 73211.43
$ws.Range("K129").Value = 16120.2
$ws.Range("L129").Value = 219634.29
$ws.Range("M129").Value = -11120.2
$ws.Range("N129").Value = -229634.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2001.25
$ws.Range("I131").Value = 1886.25
$ws.Range("J131").Value = 2058.75
$ws.Range("K131").Value = 5658.75
$ws.Range("L131").Value = 6176.25
$ws.Range("M131").Value = -618.75
$ws.Range("N131").Value = -16256.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 643.3333
$ws.Range("I140").Value = 643.3333
$ws.Range("K140").Value = 1929.9999
$ws.Range("M140").Value = 3250.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5600
$ws.Range("J113").Value = 7500
$ws.Range("N113").Value = -11840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 66332.336
$ws.Range("J112").Value = 66332.336
$ws.Range("L112").Value = 66332.336
$ws.Range("N112").Value = -69286.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1708.7368
$ws.Range("I136").Value = 1290
$ws.Range("J136").Value = 2736.5454
$ws.Range("K136").Value = 3870
$ws.Range("L136").Value = 8209.6362
$ws.Range("M136").Value = -1320
$ws.Range("N136").Value = -13309.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16032.417
$ws.Range("I132").Value = 1786.3125
$ws.Range("K132").Value = 5358.9375
$ws.Range("M132").Value = -2828.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2379.0527
$ws.Range("I136").Value = 933.5
$ws.Range("J136").Value = 4857.143
$ws.Range("K136").Value = 2800.5
$ws.Range("L136").Value = 14571.429
$ws.Range("M136").Value = -250.5
$ws.Range("N136").Value = -19671.429
